# "Drop unit test to LvL2"
# The "Application must feature unit testing" task is removed from the
# Level 1 "must have" block (was row 8) and re-appended as a new task
# at the bottom of the Level 2 block (new row 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove row 8 ("Application must feature unit testing") and shift the
#    remaining Level-1 "must" rows (9-13) up by one.
$ws.Rows("8:8").Delete(-4162)   # xlShiftUp

# 2) The delete above also pulled the blank spacer rows and the whole
#    "Level 2" block up by one row. Re-insert a spacer row so "Level 2"
#    lands back on row 16, matching the rest of the sheet layout.
$ws.Rows("13:13").Insert(-4121) # xlShiftDown

# 3) Row 13 is now a stray duplicate of the database-connection-module row
#    (copied down by the insert). Clear it back to an empty, unstyled row -
#    only the label columns (A/C/E) remain as blank placeholder cells.
$ws.Range("A13").ClearContents()
$ws.Range("A13").NumberFormat = "General"
$ws.Range("C13").NumberFormat = $ws.Range("C13").NumberFormat
$ws.Range("E13").NumberFormat = $ws.Range("E13").NumberFormat
$ws.Range("B13").Clear()
$ws.Range("D13").Clear()

# 4) Append the unit-testing task to the end of the Level 2 list (row 22).
$ws.Range("A22").Value = "Application must feature unit testing"
$ws.Range("B22").Value = "X"

# 5) Match the author's final selection/cursor position.
$ws.Range("A13").Select()
